# Add season record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from the existing last header cell (AC1)
# onto the new header cells so they match the rest of the header row.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Set the new header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record for every data row (2 through 56).
$ws.Range("AD2:AD56").Value = 78
$ws.Range("AE2:AE56").Value = 84
$ws.Range("AF2:AF56").Value = 0

Write-Host "Season record columns added"
